# The "2024" worksheet (second sheet in the workbook) has a new entry
# logged at the top of the September ("R"/"S" columns) activity list.
# This pushes all the existing September..August rows down by one row
# (rows 29-52 become rows 30-53), growing the used range from A1:Y52 to
# A1:Y53, and then fills in the new row 29 with the latest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new blank row at row 29, shifting rows 29:52 down to 30:53.
$ws.Rows.Item(29).Insert()

# Populate the freshly inserted row with the newest log entry.
$ws.Range("R29").Value = "axis"
$ws.Range("S29").Value = "2024-09-04 07:02:13"
